$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 3 (it duplicated row 2's data set and is no longer needed once
# size becomes a single dynamic dropdown-driven row)
$ws.Rows(3).Delete()

# Update row 2 with the new article's data
$ws.Range("B2").Value = "C4"          # Category: C1 -> C4
$ws.Range("C2").Value = "SC3"         # Subcategory: SC1 -> SC3
$ws.Range("G2").Value = 37            # Description: 42 -> 37
$ws.Range("I2").Value = "L.GREY"      # Color: "L. GREY" -> "L.GREY"

# Size / Brand / ItemMRP must stay text (shared-string) cells, not numbers,
# so force the cell format to Text before writing, then restore the style
# so no stray formatting lingers on the cell.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "36"          # Size: 38 -> 36
$ws.Range("J2").Style = "Normal"

$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = "4444"        # Brand: 3333 -> 4444
$ws.Range("L2").Style = "Normal"

$ws.Range("N2").Value = "Supplier4"   # Supplier: new column populated

$ws.Range("R2").NumberFormat = "@"
$ws.Range("R2").Value = "2195.00"     # ItemMRP: 2198 -> 2195.00
$ws.Range("R2").Style = "Normal"

$ws.Range("T2").Value = 1             # Quantity: 5 -> 1
